# Apply the data/formula updates to Sheet1.
#
# Background (reverse-engineered from the data):
#   Column H ("range") = (High - Low) / divisor.  The divisor is changing
#   from 10 to 2.5 (i.e. H values become 4x larger).
#   Column I ("target") = Open[r] + H[r-1]               (row 2 has no target)
#   Column J ("ror")    = IF(High[r] >= I[r], Close[r]/I[r], 1)   (row 2 -> 1)
#   Column K ("hpr")    = running cumulative product of J
#   Column L ("dd")     = (runningMax(K) - K) / runningMax(K) * 100
#
# In addition, the last row (201) got refreshed market data for
# Close/Volume/Value (columns E/F/G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 201

# Column indices (1-based): A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8 I=9 J=10 K=11 L=12

# 1) Refresh the last row's raw market data (Close / Volume / Value).
$ws.Cells.Item($lastRow, 5).Value2 = 62850
$ws.Cells.Item($lastRow, 6).Value2 = 5492536.66695731
$ws.Cells.Item($lastRow, 7).Value2 = 340323142299.8795

# 2) Recompute column H ("range") for every data row using the new divisor.
for ($r = 2; $r -le $lastRow; $r++) {
    $high = $ws.Cells.Item($r, 3).Value2
    $low  = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 8).Value2 = ($high - $low) / 2.5
}

# 3) Recompute columns I, J, K, L in a single forward pass, since each row
#    depends on the prior row's H/K values.
$maxK = 0
$prevK = 0
$prevH = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $open  = $ws.Cells.Item($r, 2).Value2
    $high  = $ws.Cells.Item($r, 3).Value2
    $close = $ws.Cells.Item($r, 5).Value2
    $h     = $ws.Cells.Item($r, 8).Value2

    if ($r -eq 2) {
        # First row has no target / prior range to reference.
        $j = 1
    } else {
        $target = $open + $prevH
        $ws.Cells.Item($r, 9).Value2 = $target

        if ($high -ge $target) {
            $j = $close / $target
        } else {
            $j = 1
        }
    }

    $ws.Cells.Item($r, 10).Value2 = $j

    if ($r -eq 2) {
        $k = $j
    } else {
        $k = $prevK * $j
    }
    $ws.Cells.Item($r, 11).Value2 = $k

    if ($k -gt $maxK) {
        $maxK = $k
    }
    $dd = ($maxK - $k) / $maxK * 100
    $ws.Cells.Item($r, 12).Value2 = $dd

    $prevK = $k
    $prevH = $h
}
